# "Generate Report for Handback" -- mark the two localization files as
# handed back (in sync with en-US) and record the handback file/datetime
# for each locale (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (B, C) for both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-locale detail sheets: Status column (C), Latest Handback DateTime
# (H), and the new Latest Target File / Latest Handback File hyperlinks
# (F, G) that mirror the existing Latest Handoff File (D) links.
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn"
        HandbackDateTime = "2016-03-21 08:55:32"
        MdBase = "https://github.com/OpenLocalizationTest/oltest/blob/b1cf16bb38b4e0d5165effcd622e2a2c10113758/e2e"
        XlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d3e312a58c07d8c318fb7eb65fcb039d54526f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
        Xlf1 = "38616b35-e4ca-474e-89ba-bc5ae8ca2baf.32950f0f26a77ad2d894dfbdfb39941120c58269.zh-cn.xlf"
        Xlf2 = "7d7de137-a926-4996-aacb-bf98e1e53abb.8ccd51b71a7090d57045dbd21e5ba5e062ed2214.zh-cn.xlf"
    },
    @{
        Sheet = "de-de"
        HandbackDateTime = "2016-03-21 08:55:38"
        MdBase = "https://github.com/OpenLocalizationTest/oltest/blob/b1cf16bb38b4e0d5165effcd622e2a2c10113758/e2e"
        XlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c744c20b4d5414fc4ae7e05a995a60cb2130699/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
        Xlf1 = "38616b35-e4ca-474e-89ba-bc5ae8ca2baf.32950f0f26a77ad2d894dfbdfb39941120c58269.de-de.xlf"
        Xlf2 = "7d7de137-a926-4996-aacb-bf98e1e53abb.8ccd51b71a7090d57045dbd21e5ba5e062ed2214.de-de.xlf"
    }
)

$md1 = "38616b35-e4ca-474e-89ba-bc5ae8ca2baf.md"
$md2 = "7d7de137-a926-4996-aacb-bf98e1e53abb.md"

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status + handback datetime for both rows.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Range("H2").Value = $locale.HandbackDateTime
    $ws.Range("H3").Value = $locale.HandbackDateTime

    # Row 2 (first file): Latest Target File (F2) + Latest Handback File (G2).
    $ws.Hyperlinks.Add($ws.Range("F2"), ($locale.MdBase + "/" + $md1), "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("G2"), ($locale.XlfBase + "/" + $locale.Xlf1), "", "", $locale.Xlf1)

    # Row 3 (second file): Latest Target File (F3) + Latest Handback File (G3).
    $ws.Hyperlinks.Add($ws.Range("F3"), ($locale.MdBase + "/" + $md2), "", "", $md2)
    $ws.Hyperlinks.Add($ws.Range("G3"), ($locale.XlfBase + "/" + $locale.Xlf2), "", "", $locale.Xlf2)
}
